$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column-header strings:
#      "<Base>_old" -> "<Base>_FV2310"
#      "<Base>_new" -> "<Base>_FV2404"
#    (the "diff" header in column K is left untouched)
# ---------------------------------------------------------------------------
$headersFV2310 = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)

$headersFV2404 = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the A1:U92 range into a native Excel table ("Table1"), without
#    letting the table-creation path inject a header dxf: stash the header
#    row's current formatting in a scratch row, clear the header format
#    (creating a table over an already-bold/filled/bordered header is what
#    triggers an automatic headerRowDxfId), add the table, then restore the
#    original header formatting from the stash in one shot.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A95:U95")

$headerRange.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(1, $tableRange, [Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false

$scratch.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, frozen, bottom-left pane
#    active) - the classic "freeze top row" recipe.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
